$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder columns: date (B) moves to the end (before transaction_id),
#     category/description/amount/transaction_id shift left.
#     Cut+Insert also carries the exact (sub-pixel) column width of the
#     date column along with it, landing on the new date column (E).
$ws.Columns("B:B").Cut() | Out-Null
$ws.Columns("F:F").Insert() | Out-Null

# --- The cut/insert above leaves stray zero-width column overrides on
#     B:D (now category/description/amount) because their physical
#     position shifted. Clear them by deleting+reinserting A:D (anchored
#     on the untouched id column) which resets the column metadata; the
#     values are rewritten explicitly right after.
$ws.Columns("A:D").Delete() | Out-Null
$ws.Columns("A:D").Insert() | Out-Null

# --- Header row: rename + reorder the column titles.
$ws.Range("A1").Value = "Income ID"
$ws.Range("B1").Value = "Category"
$ws.Range("C1").Value = "Description"
$ws.Range("D1").Value = "Amount"
$ws.Range("E1").Value = "Date"
$ws.Range("F1").Value = "Transaction ID"

# --- Data rows: id, category, description, amount (date + transaction_id
#     already correctly in place from the column move above).
$data = @(
    @(1,  "Salary",      "FT Job",                   2500),
    @(2,  "Freelance",   "Graphic Design for Apple",  1000),
    @(3,  "Investments", "Apple Stock",               102.2),
    @(4,  "Salary",      "FT Job",                    2500),
    @(5,  "Gifts",       "CNY Angpaos",                288),
    @(6,  "Investments", "Apple Stock",                58.7),
    @(7,  "Salary",      "FT Job",                    2500),
    @(8,  "Investments", "Apple Stock",               17.010000000000002),
    @(9,  "Salary",      "FT Job",                    2500),
    @(10, "Freelance",   "Graphic Design for SAF",    1000),
    @(11, "Investments", "Apple Stock",               45.67)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# --- Selection, matching the committed state.
$ws.Range("F2").Select() | Out-Null
